$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.437.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6293"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.117.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07742"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.035"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6809"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001064"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.46"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.188"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.520.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.456"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.437"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.390"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05612"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.134"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.590"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01804"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.231.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.739"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.449"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9100"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.204"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4026"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1156"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.029"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000116"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.37%  "

$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05709"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
